$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Skipped", "No Easy Apply", "19-01-26", "21:32", "https://www.linkedin.com/jobs/search/?currentJobId=4323316446&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland"),
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Skipped", "No Easy Apply", "19-01-26", "21:32", "https://www.linkedin.com/jobs/search/?currentJobId=4362486265&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland"),
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Skipped", "No Easy Apply", "19-01-26", "21:38", "https://www.linkedin.com/jobs/search/?currentJobId=4323316446&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland"),
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Skipped", "No Easy Apply", "19-01-26", "21:39", "https://www.linkedin.com/jobs/search/?currentJobId=4362486265&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland"),
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Skipped", "No Easy Apply", "19-01-26", "21:40", "https://www.linkedin.com/jobs/search/?currentJobId=4323316446&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland"),
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Skipped", "No Easy Apply", "19-01-26", "21:40", "https://www.linkedin.com/jobs/search/?currentJobId=4362486265&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland"),
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Skipped", "No Easy Apply", "19-01-26", "21:43", "https://www.linkedin.com/jobs/search/?currentJobId=4323316446&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland"),
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Skipped", "No Easy Apply", "19-01-26", "21:44", "https://www.linkedin.com/jobs/search/?currentJobId=4362486265&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland"),
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Applied", "Success", "19-01-26", "22:02", "https://www.linkedin.com/jobs/search/?currentJobId=4354434581&f_AL=true&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland&origin=JOB_SEARCH_PAGE_JOB_FILTER&spellCorrectionEnabled=true")
)

$startRow = 6
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
